# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the
# 2c9f4975-a5f1-4a10-8518-c03d0e449099 file's xliff handoff/handback rows,
# reflecting a freshly regenerated handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 corresponds to 2c9f4975-...-c03d0e449099.md
# Column G = "Latest HO Xliff Generate Date"
$wsOverview.Range("G3").Value = "2016-09-07 08:44:02"

# zh-cn sheet: row 3 corresponds to the 2c9f4975-... zh-cn xliff
# Column H = "Correspond Handoff Datetime"
# Column K = "Correspond Handback DateTime"
$wsZhCn.Range("H3").Value = "2016-09-07 08:43:50"
$wsZhCn.Range("K3").Value = "2016-09-07 08:44:59"

# de-de sheet: row 3 corresponds to the 2c9f4975-... de-de xliff
# Column H = "Correspond Handoff Datetime"
$wsDeDe.Range("H3").Value = "2016-09-07 08:45:35"
